$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was recorded for this market/product, which
# needs to be inserted as a new data row right after the existing row for
# this date range (pushing the old row 92 and everything below it down by
# one row). Excel's native row-insert handles the shift of all subsequent
# rows (and their formatting) automatically.
$ws.Rows("92:92").Insert()

# Populate the newly-inserted row 92 with the new observation.
$ws.Range("A92").Value = 7
$ws.Range("B92").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C92").Value = "Ñuble"
$ws.Range("D92").Value = "2023-08-16"
$ws.Range("E92").Value = 16
$ws.Range("F92").Value = 100112001
$ws.Range("G92").Value = "Berenjena"
$ws.Range("H92").Value = "Sin especificar"
$ws.Range("I92").Value = "Primera"
$ws.Range("J92").Value = 50
$ws.Range("K92").Value = 9000
$ws.Range("L92").Value = 9000
$ws.Range("M92").Value = 9000
$ws.Range("N92").Value = "$/caja 60 unidades"
$ws.Range("O92").Value = "Región de Arica y Parinacota"
$ws.Range("P92").Value = 150
$ws.Range("Q92").Value = 60
$ws.Range("R92").Value = "Hortaliza"
